$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -1.182007550723355
$ws1.Range("C2").Value = -1.012740346017771
$ws1.Range("B3").Value = -0.3861856897723882
$ws1.Range("C3").Value = -0.8808188739039313
$ws1.Range("B4").Value = -0.3732814333320856
$ws1.Range("C4").Value = 0.09855202579779687

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -1.054558373124875
$ws2.Range("C2").Value = -0.1716155305906732
$ws2.Range("B3").Value = 0.2318440616085357
$ws2.Range("C3").Value = 0.878025215190953
$ws2.Range("B4").Value = -0.3487201350163849
$ws2.Range("C4").Value = 0.08863323165777205
